$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "24/10/2025"
$ws.Range("B14").Value = "Lok. Sofia"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = "Lok. Plovdiv"
$ws.Range("F14").Value = "D"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 1
$ws.Range("M14").Value = 8
$ws.Range("N14").Value = 13
$ws.Range("O14").Value = 4
$ws.Range("P14").Value = 4
